# Stage 1: update companies data
# Swap the data of several rows in the relevant_companies sheet back to
# reflect the corrected company ordering / details.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Company Number (B) and SIC Codes (I) columns must stay text
# (some values are purely numeric, e.g. "16473813"), so force text
# formatting before writing values to avoid Excel auto-converting them
# to numbers.
$ws.Range("B2:B11").NumberFormat = "@"
$ws.Range("I2:I11").NumberFormat = "@"

# Row 2 <-> Row 4 swap (GANDER INVESTMENTS LTD <-> BDL NORTH PARTNERS LTD)
# (J2/K2 and J4/K4 are blank both before and after the edit, so they are
# left untouched.)
$ws.Range("A2").Value = "BDL NORTH PARTNERS LTD"
$ws.Range("B2").Value = "16473813"
$ws.Range("H2").Value = "Partners"
$ws.Range("I2").Value = "68320"

$ws.Range("A4").Value = "GANDER INVESTMENTS LTD"
$ws.Range("B4").Value = "16473515"
$ws.Range("H4").Value = "Investments"
$ws.Range("I4").Value = "68100,68209"

# Row 5 now holds INTERCONTINENTAL HOLDING COMPANY LIMITED's data (was row 9)
$ws.Range("A5").Value = "INTERCONTINENTAL HOLDING COMPANY LIMITED"
$ws.Range("B5").Value = "16473418"
$ws.Range("H5").Value = "Other"
$ws.Range("I5").Value = "64209"
$ws.Range("J5").Value = "Activities of other holding companies n.e.c."
$ws.Range("K5").Value = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles."

# Row 6 now holds AJ INVESTMENT AND CONSULTANCY LTD's data (was row 11)
$ws.Range("A6").Value = "AJ INVESTMENT AND CONSULTANCY LTD"
$ws.Range("B6").Value = "16473328"
$ws.Range("H6").Value = "Investments"
$ws.Range("I6").Value = "64306,70229"
$ws.Range("J6").Value = "Activities of real estate investment trusts"
$ws.Range("K6").Value = "UK-regulated REIT companies."

# Row 7 now holds GAUNT CAPITAL LTD's data (was row 10)
$ws.Range("A7").Value = "GAUNT CAPITAL LTD"
$ws.Range("B7").Value = "16473262"
$ws.Range("H7").Value = "Capital"
$ws.Range("I7").Value = "64209"
$ws.Range("J7").Value = "Activities of other holding companies n.e.c."
$ws.Range("K7").Value = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles."

# Row 9 now holds THE DISLEY GROUP LTD's data (was row 7)
$ws.Range("A9").Value = "THE DISLEY GROUP LTD"
$ws.Range("B9").Value = "16473398"
$ws.Range("H9").Value = "Other"
$ws.Range("I9").Value = "64209"
$ws.Range("J9").Value = "Activities of other holding companies n.e.c."
$ws.Range("K9").Value = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles."

# Row 10 now holds MARMIMI HOLDING LIMITED's data (was row 6)
$ws.Range("A10").Value = "MARMIMI HOLDING LIMITED"
$ws.Range("B10").Value = "16473234"
$ws.Range("H10").Value = "Other"
$ws.Range("I10").Value = "64209"
$ws.Range("J10").Value = "Activities of other holding companies n.e.c."
$ws.Range("K10").Value = "Catch-all SPV: protected cells, cell companies, bespoke feeder vehicles."

# Row 11 now holds BRIDGEWICK PARTNERS LIMITED's data (was row 5)
$ws.Range("A11").Value = "BRIDGEWICK PARTNERS LIMITED"
$ws.Range("B11").Value = "16473142"
$ws.Range("H11").Value = "Partners"
$ws.Range("I11").Value = "64999"
$ws.Range("J11").Value = "Financial intermediation not elsewhere classified"
$ws.Range("K11").Value = "Catch-all credit-oriented SPVs for novel lending structures."
